$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1279.236
$ws.Range("I15").Value = 1279.236
$ws.Range("K15").Value = 3837.708000000001
$ws.Range("M15").Value = -3668.708000000001
$ws.Range("H69").Value = 4588.8887
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 4787.5
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 14362.5
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -16110.5
$ws.Range("H72").Value = 4588.8887
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 4787.5
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 43087.5
$ws.Range("M72").Value = -22632
$ws.Range("N72").Value = -51823.5
$ws.Range("H98").Value = 1019.2632
$ws.Range("I98").Value = 925.1892
$ws.Range("K98").Value = 925.1892
$ws.Range("M98").Value = 572.8108
$ws.Range("H122").Value = 1019.2632
$ws.Range("I122").Value = 925.1892
$ws.Range("K122").Value = 2775.5676
$ws.Range("M122").Value = -325.5676000000003
$ws.Range("H127").Value = 4406.6
$ws.Range("I127").Value = 5058.25
$ws.Range("J127").Value = 1800
$ws.Range("K127").Value = 15174.75
$ws.Range("L127").Value = 5400
$ws.Range("M127").Value = -10214.75
$ws.Range("N127").Value = -15320
$ws.Range("H137").Value = 68590.336
$ws.Range("I137").Value = 1818.375
$ws.Range("J137").Value = 144901.14
$ws.Range("K137").Value = 5455.125
$ws.Range("L137").Value = 434703.42
$ws.Range("M137").Value = -2905.125
$ws.Range("N137").Value = -439803.42

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8082.65
$ws.Range("I32").Value = 6218.4717
$ws.Range("J32").Value = 23165.545
$ws.Range("K32").Value = 6218.4717
$ws.Range("L32").Value = 23165.545
$ws.Range("M32").Value = -5931.4717
$ws.Range("N32").Value = -23739.545
$ws.Range("H45").Value = 11251574
$ws.Range("I45").Value = 45000600
$ws.Range("K45").Value = 45000600
$ws.Range("M45").Value = -45000223

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2000000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 2000000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("H94").Value = 172.875
$ws.Range("I94").Value = 177.57143
$ws.Range("K94").Value = 177.57143
$ws.Range("M94").Value = 273.42857
$ws.Range("H99").Value = 2405
$ws.Range("I99").Value = 2405
$ws.Range("K99").Value = 2405
$ws.Range("M99").Value = -907
$ws.Range("H105").Value = 2322.7778
$ws.Range("I105").Value = 2047.6471
$ws.Range("K105").Value = 2047.6471
$ws.Range("M105").Value = -300.6470999999999
$ws.Range("H134").Value = 4476.5454
$ws.Range("I134").Value = 4465.778
$ws.Range("K134").Value = 13397.334
$ws.Range("M134").Value = -10862.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 61
$ws.Range("J7").Value = 50
$ws.Range("L7").Value = 50
$ws.Range("N7").Value = -276
$ws.Range("H22").Value = 1548.909
$ws.Range("I22").Value = 760
$ws.Range("K22").Value = 760
$ws.Range("M22").Value = -410
$ws.Range("H107").Value = 1912.5217
$ws.Range("J107").Value = 2399.7144
$ws.Range("L107").Value = 2399.7144
$ws.Range("N107").Value = -6239.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 83392860
$ws.Range("J50").Value = 166668800
$ws.Range("L50").Value = 500006400
$ws.Range("N50").Value = -500007362
$ws.Range("H53").Value = 83392860
$ws.Range("J53").Value = 166668800
$ws.Range("L53").Value = 500006400
$ws.Range("N53").Value = -500007362
$ws.Range("H55").Value = 5000
$ws.Range("J55").Value = 5000
$ws.Range("L55").Value = 15000
$ws.Range("N55").Value = -15354
$ws.Range("H96").Value = 5997.778
$ws.Range("J96").Value = 5997.778
$ws.Range("L96").Value = 17993.334
$ws.Range("N96").Value = -22111.334
$ws.Range("H97").Value = 976.1818
$ws.Range("I97").Value = 256.66666
$ws.Range("K97").Value = 769.9999799999999
$ws.Range("M97").Value = -273.9999799999999
$ws.Range("H107").Value = 658
$ws.Range("I107").Value = 492.7857
$ws.Range("J107").Value = 812.2
$ws.Range("K107").Value = 1478.3571
$ws.Range("L107").Value = 2436.6
$ws.Range("M107").Value = 441.6428999999998
$ws.Range("N107").Value = -6276.6
$ws.Range("H113").Value = 40284.57
$ws.Range("J113").Value = 1149.3158
$ws.Range("L113").Value = 3447.9474
$ws.Range("N113").Value = -7787.9474

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 900.55554
$ws.Range("I122").Value = 700.625
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 2101.875
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = 348.125
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 2573470.5
$ws.Range("I126").Value = 3271009.5
$ws.Range("K126").Value = 9813028.5
$ws.Range("M126").Value = -9810558.5
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H132").Value = 877996.8
$ws.Range("J132").Value = 2881
$ws.Range("L132").Value = 8643
$ws.Range("N132").Value = -13703

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1654
$ws.Range("I93").Value = 1080.6666
$ws.Range("K93").Value = 1080.6666
$ws.Range("M93").Value = 167.3334
$ws.Range("H100").Value = 1866
$ws.Range("I100").Value = 1399
$ws.Range("K100").Value = 1399
$ws.Range("M100").Value = -858
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H122").Value = 2808.7778
$ws.Range("I122").Value = 2611.2856
$ws.Range("K122").Value = 7833.8568
$ws.Range("M122").Value = -5383.8568
$ws.Range("H132").Value = 2660.3794
$ws.Range("I132").Value = 1902.4445
$ws.Range("J132").Value = 3320.516
$ws.Range("K132").Value = 5707.333500000001
$ws.Range("L132").Value = 9961.548
$ws.Range("M132").Value = -3177.333500000001
$ws.Range("N132").Value = -15021.548
$ws.Range("H136").Value = 1407.3846
$ws.Range("I136").Value = 929.2
$ws.Range("J136").Value = 3001.3333
$ws.Range("K136").Value = 2787.6
$ws.Range("L136").Value = 9003.999899999999
$ws.Range("M136").Value = -237.6000000000004
$ws.Range("N136").Value = -14103.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 47098
$ws.Range("J70").Value = 47098
$ws.Range("L70").Value = 47098
$ws.Range("N70").Value = -47728
$ws.Range("H73").Value = 47098
$ws.Range("J73").Value = 47098
$ws.Range("L73").Value = 47098
$ws.Range("N73").Value = -49282
$ws.Range("H107").Value = 612.32
$ws.Range("I107").Value = 495.45
$ws.Range("J107").Value = 1079.8
$ws.Range("K107").Value = 1486.35
$ws.Range("L107").Value = 3239.4
$ws.Range("M107").Value = 433.6500000000001
$ws.Range("N107").Value = -7079.4
$ws.Range("H132").Value = 1703.375
$ws.Range("I132").Value = 1537.762
$ws.Range("K132").Value = 4613.286
$ws.Range("M132").Value = -2083.286
$ws.Range("H136").Value = 10894301
$ws.Range("I136").Value = 16835758
$ws.Range("K136").Value = 50507274
$ws.Range("M136").Value = -50504724
$ws.Range("H139").Value = 69376.92
$ws.Range("J139").Value = 69376.92
$ws.Range("L139").Value = 69376.92
$ws.Range("N139").Value = -79656.92

